$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NTTYY")

$ws.Range("D8").Value = 106668300
$ws.Range("E8").Value = 102974800
$ws.Range("F8").Value = 104330600
$ws.Range("G8").Value = 100301700
$ws.Range("H8").Value = 98763600
$ws.Range("I8").Value = 96734700
$ws.Range("J8").Value = 94986600
$ws.Range("D9").Value = 51848300
$ws.Range("E9").Value = 49976000
$ws.Range("F9").Value = 50859400
$ws.Range("G9").Value = 47768200
$ws.Range("H9").Value = 44207300
$ws.Range("I9").Value = 41314400
$ws.Range("J9").Value = 39567500
$ws.Range("D10").Value = 54820000
$ws.Range("E10").Value = 52998800
$ws.Range("F10").Value = 53471200
$ws.Range("G10").Value = 52533500
$ws.Range("H10").Value = 54556200
$ws.Range("I10").Value = 55420200
$ws.Range("J10").Value = 55419000
$ws.Range("I12").Value = 2433500
$ws.Range("J12").Value = 2421400
$ws.Range("D14").Value = 1466000
$ws.Range("E14").Value = 667600
$ws.Range("F14").Value = 295800
$ws.Range("G14").Value = 381500
$ws.Range("H14").Value = 107800
$ws.Range("I14").Value = 332300
$ws.Range("J14").Value = 135800
$ws.Range("D15").Value = 12108400
$ws.Range("E15").Value = 13218600
$ws.Range("F15").Value = 15967600
$ws.Range("G15").Value = 16525100
$ws.Range("H15").Value = 16997800
$ws.Range("I15").Value = 34338300
$ws.Range("J15").Value = 17272700
$ws.Range("D17").Value = 91817000
$ws.Range("E17").Value = 89055100
$ws.Range("F17").Value = 92143300
$ws.Range("G17").Value = 90497200
$ws.Range("H17").Value = 87792100
$ws.Range("I17").Value = 85868900
$ws.Range("J17").Value = 83930900
$ws.Range("D18").Value = 14851300
$ws.Range("E18").Value = 13919700
$ws.Range("F18").Value = 12187300
$ws.Range("G18").Value = 9804500
$ws.Range("H18").Value = 10971400
$ws.Range("I18").Value = 10865800
$ws.Range("J18").Value = 11055600
$ws.Range("D20").Value = 1310500
$ws.Range("E20").Value = 232700
$ws.Range("F20").Value = 205900
$ws.Range("G20").Value = 235800
$ws.Range("H20").Value = 1159200
$ws.Range("I20").Value = 452200
$ws.Range("J20").Value = 657100
$ws.Range("D21").Value = 28283600
$ws.Range("E21").Value = 27385600
$ws.Range("F21").Value = 28378400
$ws.Range("G21").Value = 26583600
$ws.Range("H21").Value = 29147200
$ws.Range("I21").Value = 28506100
$ws.Range("J21").Value = 29054100
$ws.Range("D22").Value = 291000
$ws.Range("E22").Value = 341400
$ws.Range("F22").Value = 376700
$ws.Range("G22").Value = 397900
$ws.Range("H22").Value = 431100
$ws.Range("I22").Value = 491200
$ws.Range("J22").Value = 509200
$ws.Range("D23").Value = 15870800
$ws.Range("E23").Value = 13811000
$ws.Range("F23").Value = 12016500
$ws.Range("G23").Value = 9642300
$ws.Range("H23").Value = 11699500
$ws.Range("I23").Value = 10826700
$ws.Range("J23").Value = 11203500
$ws.Range("D24").Value = 4898500
$ws.Range("E24").Value = 4234100
$ws.Range("F24").Value = 3207600
$ws.Range("G24").Value = 3592000
$ws.Range("H24").Value = 4398400
$ws.Range("I24").Value = 4284500
$ws.Range("J24").Value = 5313600
$ws.Range("D26").Value = 10972400
$ws.Range("E26").Value = 9577000
$ws.Range("F26").Value = 8808900
$ws.Range("G26").Value = 6050300
$ws.Range("H26").Value = 7301100
$ws.Range("I26").Value = 6542200
$ws.Range("J26").Value = 5889900
$ws.Range("D27").Value = 8223600
$ws.Range("E27").Value = 7233200
$ws.Range("F27").Value = 6669200
$ws.Range("G27").Value = 4683300
$ws.Range("H27").Value = 5292700
$ws.Range("I27").Value = 4718300
$ws.Range("J27").Value = 4228000
$ws.Range("D32").Value = -1310500
$ws.Range("E32").Value = -232700
$ws.Range("F32").Value = -205900
$ws.Range("G32").Value = -235800
$ws.Range("H32").Value = -1159200
$ws.Range("I32").Value = -452200
$ws.Range("J32").Value = -657100
$ws.Range("D33").Value = 8223600
$ws.Range("E33").Value = 7233200
$ws.Range("F33").Value = 6669200
$ws.Range("G33").Value = 4683300
$ws.Range("H33").Value = 5292700
$ws.Range("I33").Value = 4718300
$ws.Range("J33").Value = 4228000
$ws.Range("D35").Value = 8223600
$ws.Range("E35").Value = 7233200
$ws.Range("F35").Value = 6669200
$ws.Range("G35").Value = 4683300
$ws.Range("H35").Value = 5292700
$ws.Range("I35").Value = 4718300
$ws.Range("J35").Value = 4228000
$ws.Range("D41").Value = 15144700
$ws.Range("E41").Value = 8363900
$ws.Range("F41").Value = 9838000
$ws.Range("G41").Value = 7676500
$ws.Range("H41").Value = 8899500
$ws.Range("I41").Value = 17382700
$ws.Range("J41").Value = 9222100
$ws.Range("D42").Value = 1401100
$ws.Range("E42").Value = 577100
$ws.Range("F42").Value = 299000
$ws.Range("G42").Value = 328500
$ws.Range("H42").Value = 352100
$ws.Range("I42").Value = 485900
$ws.Range("J42").Value = 2774600
$ws.Range("D43").Value = 68781300
$ws.Range("E43").Value = 28532300
$ws.Range("F43").Value = 28576100
$ws.Range("G43").Value = 27371600
$ws.Range("H43").Value = 25378300
$ws.Range("I43").Value = 24773200
$ws.Range("J43").Value = 22752800
$ws.Range("D44").Value = 6759800
$ws.Range("E44").Value = 3303000
$ws.Range("F44").Value = 3747800
$ws.Range("G44").Value = 3530300
$ws.Range("H44").Value = 3754400
$ws.Range("I44").Value = 3170500
$ws.Range("J44").Value = 2977500
$ws.Range("D45").Value = 9493000
$ws.Range("E45").Value = 7247900
$ws.Range("F45").Value = 6599000
$ws.Range("G45").Value = 5906300
$ws.Range("H45").Value = 5559200
$ws.Range("I45").Value = 5089400
$ws.Range("J45").Value = 4868800
$ws.Range("D46").Value = 53057200
$ws.Range("E46").Value = 48024300
$ws.Range("F46").Value = 49059900
$ws.Range("G46").Value = 44813300
$ws.Range("H46").Value = 43943500
$ws.Range("I46").Value = 42210400
$ws.Range("J46").Value = 42595900
$ws.Range("D47").Value = 23831700
$ws.Range("E47").Value = 8858200
$ws.Range("F47").Value = 8949300
$ws.Range("G47").Value = 9562800
$ws.Range("H47").Value = 8401800
$ws.Range("I47").Value = 8218300
$ws.Range("J47").Value = 7580300
$ws.Range("D48").Value = 177506000
$ws.Range("E48").Value = 87859900
$ws.Range("F48").Value = 86349400
$ws.Range("G48").Value = 88605300
$ws.Range("H48").Value = 88950800
$ws.Range("I48").Value = 88381300
$ws.Range("J48").Value = 88649500
$ws.Range("D49").Value = 48621500
$ws.Range("E49").Value = 26921600
$ws.Range("F49").Value = 25616300
$ws.Range("G49").Value = 25742900
$ws.Range("H49").Value = 25291600
$ws.Range("I49").Value = 36721600
$ws.Range("J49").Value = 21512800
$ws.Range("D52").Value = 30822500
$ws.Range("E52").Value = 20439000
$ws.Range("F52").Value = 20189900
$ws.Range("G52").Value = 18425600
$ws.Range("H52").Value = 16788300
$ws.Range("I52").Value = 15827400
$ws.Range("J52").Value = 14944400
$ws.Range("D54").Value = 194735000
$ws.Range("E54").Value = 192103000
$ws.Range("F54").Value = 190165000
$ws.Range("G54").Value = 187150000
$ws.Range("H54").Value = 183376000
$ws.Range("I54").Value = 176724000
$ws.Range("J54").Value = 175283000
$ws.Range("D57").Value = 16378000
$ws.Range("E57").Value = 14581500
$ws.Range("F57").Value = 14218100
$ws.Range("G57").Value = 14279300
$ws.Range("H57").Value = 13923900
$ws.Range("I57").Value = 12987300
$ws.Range("J57").Value = 13402600
$ws.Range("D58").Value = 15416900
$ws.Range("E58").Value = 8348800
$ws.Range("F58").Value = 5615100
$ws.Range("G58").Value = 6520600
$ws.Range("H58").Value = 6434000
$ws.Range("I58").Value = 13563900
$ws.Range("J58").Value = 6863000
$ws.Range("D59").Value = 19934300
$ws.Range("E59").Value = 14417500
$ws.Range("F59").Value = 13986600
$ws.Range("G59").Value = 12844100
$ws.Range("H59").Value = 12880700
$ws.Range("I59").Value = 11441200
$ws.Range("J59").Value = 11322400
$ws.Range("D60").Value = 41786800
$ws.Range("E60").Value = 37347800
$ws.Range("F60").Value = 33819800
$ws.Range("G60").Value = 33644100
$ws.Range("H60").Value = 33238500
$ws.Range("I60").Value = 31634500
$ws.Range("J60").Value = 31588000
$ws.Range("D61").Value = 26702800
$ws.Range("E61").Value = 28874200
$ws.Range("F61").Value = 32307500
$ws.Range("G61").Value = 33657800
$ws.Range("H61").Value = 31817400
$ws.Range("I61").Value = 29568800
$ws.Range("J61").Value = 32062500
$ws.Range("D62").Value = 42075800
$ws.Range("E62").Value = 21391500
$ws.Range("F62").Value = 22019600
$ws.Range("G62").Value = 19702200
$ws.Range("H62").Value = 19325600
$ws.Range("I62").Value = 20401400
$ws.Range("J62").Value = 20800800
$ws.Range("D66").Value = 112807000
$ws.Range("E66").Value = 110269000
$ws.Range("F66").Value = 110307000
$ws.Range("G66").Value = 108666000
$ws.Range("H66").Value = 106433000
$ws.Range("I66").Value = 102311000
$ws.Range("J66").Value = 104024000
$ws.Range("D72").Value = 55490700
$ws.Range("E72").Value = 50860400
$ws.Range("F72").Value = 45871100
$ws.Range("G72").Value = 46345000
$ws.Range("H72").Value = 43467600
$ws.Range("I72").Value = 47254500
$ws.Range("J72").Value = 44194300
$ws.Range("D76").Value = 81927300
$ws.Range("E76").Value = 81834400
$ws.Range("F76").Value = 79857600
$ws.Range("G76").Value = 78484000
$ws.Range("H76").Value = 76942600
$ws.Range("I76").Value = 74412200
$ws.Range("J76").Value = 71258600
$ws.Range("D81").Value = 8223600
$ws.Range("E81").Value = 7233200
$ws.Range("F81").Value = 6669200
$ws.Range("G81").Value = 4683300
$ws.Range("H81").Value = 5292700
$ws.Range("I81").Value = 4718300
$ws.Range("J81").Value = 4228000
$ws.Range("D83").Value = 12108400
$ws.Range("E83").Value = 13218600
$ws.Range("F83").Value = 15967600
$ws.Range("G83").Value = 16525100
$ws.Range("H83").Value = 16997800
$ws.Range("I83").Value = 17169200
$ws.Range("J83").Value = 17322200
$ws.Range("D89").Value = 23789900
$ws.Range("E89").Value = 26349100
$ws.Range("F89").Value = 24533400
$ws.Range("G89").Value = 21622000
$ws.Range("H89").Value = 24660300
$ws.Range("I89").Value = 22181300
$ws.Range("J89").Value = 22675000
$ws.Range("D91").Value = -11830900
$ws.Range("E91").Value = -11767300
$ws.Range("F91").Value = -11441200
$ws.Range("G91").Value = -13062000
$ws.Range("H91").Value = -13439300
$ws.Range("I91").Value = -13904600
$ws.Range("J91").Value = -12611600
$ws.Range("D94").Value = -16649700
$ws.Range("E94").Value = -18887400
$ws.Range("F94").Value = -15908400
$ws.Range("G94").Value = -16892000
$ws.Range("H94").Value = -19045500
$ws.Range("I94").Value = -16057300
$ws.Range("J94").Value = -17820100
$ws.Range("D96").Value = -2453500
$ws.Range("E96").Value = -2241900
$ws.Range("F96").Value = -1809600
$ws.Range("G96").Value = -1805900
$ws.Range("H96").Value = -1683000
$ws.Range("I96").Value = -1658000
$ws.Range("J96").Value = -1518500
$ws.Range("D100").Value = -8422200
$ws.Range("E100").Value = -8872900
$ws.Range("F100").Value = -6396500
$ws.Range("G100").Value = -6129200
$ws.Range("H100").Value = -5626900
$ws.Range("I100").Value = -6736400
$ws.Range("J100").Value = -8570400
$ws.Range("D101").Value = -28100
$ws.Range("E101").Value = -62900
$ws.Range("F101").Value = -67100
$ws.Range("G101").Value = 176200
$ws.Range("H101").Value = 220300
$ws.Range("I101").Value = 81700
$ws.Range("J101").Value = -36300
$ws.Range("D102").Value = -1310000
$ws.Range("E102").Value = -1474100
$ws.Range("F102").Value = 2161500
$ws.Range("G102").Value = -1223000
$ws.Range("H102").Value = 208200
$ws.Range("I102").Value = -530700
$ws.Range("J102").Value = -3751700
